$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F6").Value = 1122
$ws.Range("F7").Value = 901
$ws.Range("F11").Value = 875
$ws.Range("F12").Value = 309
$ws.Range("F14").Value = 518
$ws.Range("F15").Value = 1364
$ws.Range("F17").Value = 1244
$ws.Range("F18").Value = 2924
$ws.Range("F19").Value = 9
$ws.Range("F20").Value = 1518
$ws.Range("F21").Value = 1285
$ws.Range("F22").Value = 746
$ws.Range("F23").Value = 210
$ws.Range("F26").Value = 1055
$ws.Range("F27").Value = 372
$ws.Range("F28").Value = 3247
$ws.Range("F29").Value = 635
$ws.Range("F31").Value = 1447

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F8").Value = 3
$ws.Range("F10").Value = 17

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 766

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 766
$ws.Range("F10").Value = 1122
$ws.Range("F11").Value = 901
$ws.Range("F18").Value = 3
$ws.Range("F21").Value = 17
$ws.Range("F22").Value = 875
$ws.Range("F23").Value = 309
$ws.Range("F25").Value = 518
$ws.Range("F26").Value = 1364
$ws.Range("F28").Value = 1244
$ws.Range("F29").Value = 2924
$ws.Range("F30").Value = 9
$ws.Range("F31").Value = 1518
$ws.Range("F32").Value = 1285
$ws.Range("F33").Value = 746
$ws.Range("F34").Value = 210
$ws.Range("F39").Value = 1055
$ws.Range("F40").Value = 372
$ws.Range("F41").Value = 3247
$ws.Range("F42").Value = 635
$ws.Range("F44").Value = 1447
